# Replicate the CSV-derived report onto a second worksheet.
# The new sheet is an exact duplicate of the first sheet (same headers,
# rows, and formatting), named "currency_filtering_test1", placed right
# after the original "currency_filtering_test" sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Copy the existing sheet so it lands immediately after itself.
$ws1.Copy($null, $ws1)

# The copy becomes the new last sheet and the active sheet; grab + rename it.
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "currency_filtering_test1"

# Restore the original sheet as the active/selected one.
$ws1.Activate()
